$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet "Training" (sheet1): a new (Input, Class) pair for person-birthdate
# ("On what date was * born?") is inserted right after row 29, pushing the
# existing rows 30..72 down by one. Seven more rows are appended at the end
# for the new "health-condition_cause" extra inputs and the new
# "place-height" class.
# ---------------------------------------------------------------------------
for ($r = 72; $r -ge 30; $r--) {
    $ws1.Cells.Item($r + 1, 1).Value2 = $ws1.Cells.Item($r, 1).Value2
    $ws1.Cells.Item($r + 1, 2).Value2 = $ws1.Cells.Item($r, 2).Value2
}
$ws1.Cells.Item(30, 1).Value2 = "On what date was * born?"
$ws1.Cells.Item(30, 2).Value2 = "person-birthdate"

$ws1.Cells.Item(73, 1).Value2 = "How do I treat *?"
$ws1.Cells.Item(73, 2).Value2 = "health-condition_cause"
$ws1.Cells.Item(74, 1).Value2 = "Why do people catch the *?"
$ws1.Cells.Item(74, 2).Value2 = "health-condition_cause"
$ws1.Cells.Item(75, 1).Value2 = "What are common causes for *?"
$ws1.Cells.Item(75, 2).Value2 = "health-condition_cause"
$ws1.Cells.Item(76, 1).Value2 = "What are common reasons for people to get *?"
$ws1.Cells.Item(76, 2).Value2 = "health-condition_cause"
$ws1.Cells.Item(77, 1).Value2 = "How does a person get *?"
$ws1.Cells.Item(77, 2).Value2 = "health-condition_cause"
$ws1.Cells.Item(78, 1).Value2 = "How tall is *?"
$ws1.Cells.Item(78, 2).Value2 = "place-height"
$ws1.Cells.Item(79, 1).Value2 = "What's the height of *?"
$ws1.Cells.Item(79, 2).Value2 = "place-height"
$ws1.Cells.Item(80, 1).Value2 = "How high is *?"
$ws1.Cells.Item(80, 2).Value2 = "place-height"

# ---------------------------------------------------------------------------
# Sheet "original" (sheet2): four new rows for the "place-height" class,
# mirroring the existing table layout/formula pattern. Grow the table by
# copying the last existing data row (keeps formatting + the CONCATENATE
# formula pattern) and then overwrite the per-row values.
# ---------------------------------------------------------------------------
$ws2.Rows.Item(78).Copy()
$ws2.Rows.Item(79).Insert()
$ws2.Rows.Item(78).Copy()
$ws2.Rows.Item(79).Insert()
$ws2.Rows.Item(78).Copy()
$ws2.Rows.Item(79).Insert()
$ws2.Rows.Item(78).Copy()
$ws2.Rows.Item(79).Insert()

$ws2.Cells.Item(79, 1).Value2 = "How tall is *?"
$ws2.Cells.Item(79, 3).Value2 = "place"
$ws2.Cells.Item(79, 4).Value2 = "height"

$ws2.Cells.Item(80, 1).Value2 = "What's the height of *?"
$ws2.Cells.Item(80, 3).Value2 = "place"
$ws2.Cells.Item(80, 4).Value2 = "height"

$ws2.Cells.Item(81, 1).Value2 = "How high is *?"
$ws2.Cells.Item(81, 3).Value2 = "place"
$ws2.Cells.Item(81, 4).Value2 = "height"

$ws2.Cells.Item(82, 1).Value2 = "What's the elevation of *?"
$ws2.Cells.Item(82, 3).Value2 = "place"
$ws2.Cells.Item(82, 4).Value2 = "height"

# ---------------------------------------------------------------------------
# Selections: the committed workbook leaves the cursor on the last row of
# each table.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A80").Select()
$ws2.Activate()
$ws2.Range("A82").Select()
